$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 7
$ws1.Range("H2").Value = 26.81
$ws1.Range("L2").Value = 1.05

# Row 3
$ws1.Range("D3").Value = 7
$ws1.Range("H3").Value = 24.43
$ws1.Range("L3").Value = 1.2

# Row 4
$ws1.Range("H4").Value = 20.29
$ws1.Range("L4").Value = 0.98

# Row 5
$ws1.Range("D5").Value = 9
$ws1.Range("H5").Value = 18.17
$ws1.Range("L5").Value = 0.82

# Row 6
$ws1.Range("H6").Value = 17.65
$ws1.Range("L6").Value = 0.96

# Row 7
$ws1.Range("H7").Value = 16.65
$ws1.Range("L7").Value = 0.97

# Row 8
$ws1.Range("H8").Value = 15.23
$ws1.Range("L8").Value = 1.17

# Row 9
$ws1.Range("D9").Value = 7
$ws1.Range("H9").Value = 17.76
$ws1.Range("L9").Value = 0.86

# Row 10
$ws1.Range("D10").Value = 7
$ws1.Range("H10").Value = 17.07
$ws1.Range("L10").Value = 1.19

# Row 11
$ws1.Range("H11").Value = 12.64
$ws1.Range("L11").Value = 1.19

# Row 12
$ws1.Range("H12").Value = 11.33
$ws1.Range("L12").Value = 0.86

# Row 13
$ws1.Range("H13").Value = 10.33
$ws1.Range("L13").Value = 1.06

# Row 14
$ws1.Range("D14").Value = 7
$ws1.Range("H14").Value = 11.75
$ws1.Range("L14").Value = 1.11

# Row 15
$ws1.Range("D15").Value = 7
$ws1.Range("H15").Value = 10.75
$ws1.Range("L15").Value = 1.07

# Row 16
$ws1.Range("H16").Value = 9.75
$ws1.Range("L16").Value = 1.12

# Row 17
$ws1.Range("H17").Value = 6.95
$ws1.Range("L17").Value = 0.89

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "141"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "70"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "34"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "7"
